$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to text format so numeric-looking values
# (e.g. "1.003", "219.59") are preserved exactly as strings,
# matching the original inline-string cell type in the workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.297.36"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.697.53"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.59"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5268"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2698"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06478"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.24"
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07475"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.563"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.682.20"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5894"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008632"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.02"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.348.51"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.997"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.72"
$ws.Range("E21").Value = "  +2.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.271"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.46"
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.705"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1246"
$ws.Range("E26").Value = "  +6.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.96"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06818"
$ws.Range("E28").Value = "  +18.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.347"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.327"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.615"
$ws.Range("E31").Value = "  +3.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.569"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.672"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.034"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6228"
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.380"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.720"
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.309"
$ws.Range("E38").Value = "  +6.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01623"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.105.44"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8775"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.016"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.07"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.838.36"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.20"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.174"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05262"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4295"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.043"
$ws.Range("E51").Value = "  +3.68%  "
